$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr>'

# Collect the target paragraphs up front (indices shift only if paragraph
# count changes, which it does not here - InsertXML just replaces content).
$softUniParas = @()
$newYorkBookmarkPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "SoftUni,42.70,23.33") {
        $softUniParas += $p
    }
    if ($t -eq "New York,40.6976701,-74.2598732") {
        if ($p.Range.WordOpenXML -like "*_GoBack*") {
            $newYorkBookmarkPara = $p
        }
    }
}

# --- "New York..." paragraph that currently hosts the _GoBack bookmark -----
# The bookmark moves away from here down to the last "SoftUni" paragraph, so
# this paragraph loses it (text/runs stay identical otherwise).
if ($newYorkBookmarkPara -ne $null) {
    $xmlNoBookmark = '<w:p ' + $wns + '>' + `
        '<w:pPr>' + $rPr + '</w:pPr>' + `
        '<w:r>' + $rPr + '<w:t>New York,40.</w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r>' + $rPr + '<w:t>6976701,-</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r>' + $rPr + '<w:t>74.2598732</w:t></w:r>' + `
        '</w:p>'
    $newYorkBookmarkPara.Range.InsertXML($xmlNoBookmark)
}

# --- First "SoftUni,42.70,23.33" paragraph (Input: list) -------------------
# Becomes "Plovdiv" + ",42.70,23.33" as two runs, no bookmark.
$xmlFirst = '<w:p ' + $wns + '>' + `
    '<w:pPr>' + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t>Plovdiv</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>,42.70,23.33</w:t></w:r>' + `
    '</w:p>'
$softUniParas[0].Range.InsertXML($xmlFirst)

# --- Second "SoftUni,42.70,23.33" paragraph (Output: list, last one) -------
# Becomes "Plovdiv" run, then the relocated _GoBack bookmark, then a
# ",42.70,23.33" run.
$xmlSecond = '<w:p ' + $wns + '>' + `
    '<w:pPr>' + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t>Plovdiv</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r>' + $rPr + '<w:t>,42.70,23.33</w:t></w:r>' + `
    '</w:p>'
$softUniParas[1].Range.InsertXML($xmlSecond)

Write-Host "Done"
